$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G18").Value = "Error autenticacion"
$ws.Range("G28").Value = "Error autenticacion"
$ws.Range("G30").ClearContents()
$ws.Range("G37").Value = "Error autenticacion"
$ws.Range("G43").Value = "Error autenticacion"
$ws.Range("G53").Value = "Error autenticacion"
$ws.Range("G62").ClearContents()
$ws.Range("G64").Value = "Error autenticacion"
